$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (old MuSCs/Resolving-Mac sending-cluster rows no longer present)
$ws.Rows("6:9").Delete()

# Row 2: ECs -> ECs (target cluster column D changes from MuSCs to ECs)
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 45.76217133333333
$ws.Range("H2").Value = 137.286514
$ws.Range("I2").Value = 0.6763939203605134
$ws.Range("J2").Value = 0.6763939203605135
$ws.Range("M2").Value = 0.0006176666666666666
$ws.Range("N2").Value = 0.001853
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.02826576782688889
$ws.Range("R2").Value = 0.254391910442
$ws.Range("S2").Value = 0.6763939203605134
$ws.Range("T2").Value = 0.6763939203605135

# Row 3: FAPs -> ECs
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "ECs"
$ws.Range("G3").Value = 9.395935333333332
$ws.Range("H3").Value = 28.187806
$ws.Range("I3").Value = 0.1388778842960613
$ws.Range("J3").Value = 0.1388778842960613
$ws.Range("M3").Value = 0.0006176666666666666
$ws.Range("N3").Value = 0.001853
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.005803556057555555
$ws.Range("R3").Value = 0.052232004518
$ws.Range("S3").Value = 0.1388778842960613
$ws.Range("T3").Value = 0.1388778842960613

# Row 4: MuSCs -> ECs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "ECs"
$ws.Range("G4").Value = 12.29750866666667
$ws.Range("H4").Value = 36.892526
$ws.Range("I4").Value = 0.1817649787009828
$ws.Range("J4").Value = 0.1817649787009828
$ws.Range("M4").Value = 0.0006176666666666666
$ws.Range("N4").Value = 0.001853
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.007595761186444444
$ws.Range("R4").Value = 0.06836185067800001
$ws.Range("S4").Value = 0.1817649787009828
$ws.Range("T4").Value = 0.1817649787009828

# Row 5: Resolving-Mac -> ECs
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 0.2004796666666666
$ws.Range("H5").Value = 0.6014389999999999
$ws.Range("I5").Value = 0.002963216642442438
$ws.Range("J5").Value = 0.002963216642442439
$ws.Range("M5").Value = 0.0006176666666666666
$ws.Range("N5").Value = 0.001853
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.0001238296074444444
$ws.Range("R5").Value = 0.001114466467
$ws.Range("S5").Value = 0.002963216642442438
$ws.Range("T5").Value = 0.002963216642442439
